$d = $word.ActiveDocument

$d.Content.Find.Execute("80+16=", $true, $false, $false, $false, $false, $true, 1, $false, "28-21=", 2) | Out-Null
$d.Content.Find.Execute("84-55=", $true, $false, $false, $false, $false, $true, 1, $false, "72-27=", 2) | Out-Null
$d.Content.Find.Execute("54+30=", $true, $false, $false, $false, $false, $true, 1, $false, "89-73=", 2) | Out-Null
$d.Content.Find.Execute("63+22=", $true, $false, $false, $false, $false, $true, 1, $false, "94-19=", 2) | Out-Null
$d.Content.Find.Execute("45+45=", $true, $false, $false, $false, $false, $true, 1, $false, "74+24=", 2) | Out-Null
$d.Content.Find.Execute("60-25=", $true, $false, $false, $false, $false, $true, 1, $false, "14+9=", 2) | Out-Null
$d.Content.Find.Execute("52-40=", $true, $false, $false, $false, $false, $true, 1, $false, "9+83=", 2) | Out-Null
$d.Content.Find.Execute("95-73=", $true, $false, $false, $false, $false, $true, 1, $false, "27+44=", 2) | Out-Null
$d.Content.Find.Execute("69-26=", $true, $false, $false, $false, $false, $true, 1, $false, "67+2=", 2) | Out-Null
$d.Content.Find.Execute("71-29=", $true, $false, $false, $false, $false, $true, 1, $false, "4+17=", 2) | Out-Null
$d.Content.Find.Execute("97-54=", $true, $false, $false, $false, $false, $true, 1, $false, "65-17=", 2) | Out-Null
$d.Content.Find.Execute("19+17=", $true, $false, $false, $false, $false, $true, 1, $false, "25-22=", 2) | Out-Null
$d.Content.Find.Execute("43+19=", $true, $false, $false, $false, $false, $true, 1, $false, "49-14=", 2) | Out-Null
$d.Content.Find.Execute("99-69=", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=", 2) | Out-Null
$d.Content.Find.Execute("91-33=", $true, $false, $false, $false, $false, $true, 1, $false, "11+9=", 2) | Out-Null
$d.Content.Find.Execute("51-32=", $true, $false, $false, $false, $false, $true, 1, $false, "90-44=", 2) | Out-Null
$d.Content.Find.Execute("10+75=", $true, $false, $false, $false, $false, $true, 1, $false, "21+68=", 2) | Out-Null
$d.Content.Find.Execute("32-28=", $true, $false, $false, $false, $false, $true, 1, $false, "2+33=", 2) | Out-Null
$d.Content.Find.Execute("11+81=", $true, $false, $false, $false, $false, $true, 1, $false, "4-3=", 2) | Out-Null
$d.Content.Find.Execute("89-76=", $true, $false, $false, $false, $false, $true, 1, $false, "99-2=", 2) | Out-Null
$d.Content.Find.Execute("11+26=", $true, $false, $false, $false, $false, $true, 1, $false, "73-50=", 2) | Out-Null
$d.Content.Find.Execute("12-10=", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=", 2) | Out-Null
$d.Content.Find.Execute("18-15=", $true, $false, $false, $false, $false, $true, 1, $false, "22+0=", 2) | Out-Null
$d.Content.Find.Execute("58-8=", $true, $false, $false, $false, $false, $true, 1, $false, "58+23=", 2) | Out-Null
$d.Content.Find.Execute("84-64=", $true, $false, $false, $false, $false, $true, 1, $false, "12+69=", 2) | Out-Null
$d.Content.Find.Execute("24-2=", $true, $false, $false, $false, $false, $true, 1, $false, "78+14=", 2) | Out-Null
$d.Content.Find.Execute("51+41=", $true, $false, $false, $false, $false, $true, 1, $false, "70-54=", 2) | Out-Null
$d.Content.Find.Execute("43-21=", $true, $false, $false, $false, $false, $true, 1, $false, "55-54=", 2) | Out-Null
$d.Content.Find.Execute("32+60=", $true, $false, $false, $false, $false, $true, 1, $false, "5+7=", 2) | Out-Null
$d.Content.Find.Execute("50-13=", $true, $false, $false, $false, $false, $true, 1, $false, "90-15=", 2) | Out-Null
$d.Content.Find.Execute("12+78=", $true, $false, $false, $false, $false, $true, 1, $false, "35-11=", 2) | Out-Null
$d.Content.Find.Execute("59-38=", $true, $false, $false, $false, $false, $true, 1, $false, "87-48=", 2) | Out-Null
$d.Content.Find.Execute("79-11=", $true, $false, $false, $false, $false, $true, 1, $false, "80-44=", 2) | Out-Null
$d.Content.Find.Execute("27+58=", $true, $false, $false, $false, $false, $true, 1, $false, "72-34=", 2) | Out-Null
$d.Content.Find.Execute("83-66=", $true, $false, $false, $false, $false, $true, 1, $false, "51+29=", 2) | Out-Null
$d.Content.Find.Execute("46-12=", $true, $false, $false, $false, $false, $true, 1, $false, "45+23=", 2) | Out-Null
$d.Content.Find.Execute("34-10=", $true, $false, $false, $false, $false, $true, 1, $false, "43-3=", 2) | Out-Null
$d.Content.Find.Execute("29-21=", $true, $false, $false, $false, $false, $true, 1, $false, "58-36=", 2) | Out-Null
$d.Content.Find.Execute("84+9=", $true, $false, $false, $false, $false, $true, 1, $false, "20+78=", 2) | Out-Null
$d.Content.Find.Execute("69+26=", $true, $false, $false, $false, $false, $true, 1, $false, "90-27=", 2) | Out-Null
$d.Content.Find.Execute("47+49=", $true, $false, $false, $false, $false, $true, 1, $false, "36+33=", 2) | Out-Null
$d.Content.Find.Execute("44+27=", $true, $false, $false, $false, $false, $true, 1, $false, "44-8=", 2) | Out-Null
$d.Content.Find.Execute("16+73=", $true, $false, $false, $false, $false, $true, 1, $false, "59-23=", 2) | Out-Null
$d.Content.Find.Execute("30+10=", $true, $false, $false, $false, $false, $true, 1, $false, "24-6=", 2) | Out-Null
$d.Content.Find.Execute("74-55=", $true, $false, $false, $false, $false, $true, 1, $false, "50+33=", 2) | Out-Null
$d.Content.Find.Execute("73-32=", $true, $false, $false, $false, $false, $true, 1, $false, "32-5=", 2) | Out-Null
$d.Content.Find.Execute("26-14=", $true, $false, $false, $false, $false, $true, 1, $false, "3+13=", 2) | Out-Null
$d.Content.Find.Execute("23+0=", $true, $false, $false, $false, $false, $true, 1, $false, "48-26=", 2) | Out-Null
$d.Content.Find.Execute("25+15=", $true, $false, $false, $false, $false, $true, 1, $false, "78-72=", 2) | Out-Null
$d.Content.Find.Execute("49-21=", $true, $false, $false, $false, $false, $true, 1, $false, "23+20=", 2) | Out-Null
$d.Content.Find.Execute("61-26=", $true, $false, $false, $false, $false, $true, 1, $false, "73-23=", 2) | Out-Null
$d.Content.Find.Execute("89-43=", $true, $false, $false, $false, $false, $true, 1, $false, "19+32=", 2) | Out-Null
$d.Content.Find.Execute("2+37=", $true, $false, $false, $false, $false, $true, 1, $false, "35+34=", 2) | Out-Null
$d.Content.Find.Execute("43+45=", $true, $false, $false, $false, $false, $true, 1, $false, "14+18=", 2) | Out-Null
$d.Content.Find.Execute("24+51=", $true, $false, $false, $false, $false, $true, 1, $false, "92-28=", 2) | Out-Null
$d.Content.Find.Execute("14+27=", $true, $false, $false, $false, $false, $true, 1, $false, "66+2=", 2) | Out-Null
$d.Content.Find.Execute("82+12=", $true, $false, $false, $false, $false, $true, 1, $false, "96-26=", 2) | Out-Null
$d.Content.Find.Execute("38-30=", $true, $false, $false, $false, $false, $true, 1, $false, "50+28=", 2) | Out-Null
$d.Content.Find.Execute("62+22=", $true, $false, $false, $false, $false, $true, 1, $false, "96-94=", 2) | Out-Null
$d.Content.Find.Execute("21+24=", $true, $false, $false, $false, $false, $true, 1, $false, "8+67=", 2) | Out-Null
$d.Content.Find.Execute("89-28=", $true, $false, $false, $false, $false, $true, 1, $false, "9+11=", 2) | Out-Null
$d.Content.Find.Execute("64+32=", $true, $false, $false, $false, $false, $true, 1, $false, "13-2=", 2) | Out-Null
$d.Content.Find.Execute("47-28=", $true, $false, $false, $false, $false, $true, 1, $false, "0+27=", 2) | Out-Null
$d.Content.Find.Execute("68+1=", $true, $false, $false, $false, $false, $true, 1, $false, "77-46=", 2) | Out-Null
$d.Content.Find.Execute("21+61=", $true, $false, $false, $false, $false, $true, 1, $false, "75-22=", 2) | Out-Null
$d.Content.Find.Execute("42+7=", $true, $false, $false, $false, $false, $true, 1, $false, "17+76=", 2) | Out-Null
$d.Content.Find.Execute("74-69=", $true, $false, $false, $false, $false, $true, 1, $false, "56-7=", 2) | Out-Null
$d.Content.Find.Execute("6+5=", $true, $false, $false, $false, $false, $true, 1, $false, "3+35=", 2) | Out-Null
$d.Content.Find.Execute("84-41=", $true, $false, $false, $false, $false, $true, 1, $false, "73-54=", 2) | Out-Null
$d.Content.Find.Execute("37+45=", $true, $false, $false, $false, $false, $true, 1, $false, "4+58=", 2) | Out-Null
$d.Content.Find.Execute("5+75=", $true, $false, $false, $false, $false, $true, 1, $false, "79-36=", 2) | Out-Null
$d.Content.Find.Execute("67+12=", $true, $false, $false, $false, $false, $true, 1, $false, "24+42=", 2) | Out-Null
$d.Content.Find.Execute("55-9=", $true, $false, $false, $false, $false, $true, 1, $false, "12+45=", 2) | Out-Null
$d.Content.Find.Execute("52+14=", $true, $false, $false, $false, $false, $true, 1, $false, "31-13=", 2) | Out-Null
$d.Content.Find.Execute("93-5=", $true, $false, $false, $false, $false, $true, 1, $false, "75+17=", 2) | Out-Null
$d.Content.Find.Execute("4+33=", $true, $false, $false, $false, $false, $true, 1, $false, "93-69=", 2) | Out-Null
$d.Content.Find.Execute("54-17=", $true, $false, $false, $false, $false, $true, 1, $false, "95-88=", 2) | Out-Null
$d.Content.Find.Execute("93-62=", $true, $false, $false, $false, $false, $true, 1, $false, "48-34=", 2) | Out-Null
$d.Content.Find.Execute("87+2=", $true, $false, $false, $false, $false, $true, 1, $false, "60-36=", 2) | Out-Null
$d.Content.Find.Execute("88-59=", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=", 2) | Out-Null
$d.Content.Find.Execute("70-10=", $true, $false, $false, $false, $false, $true, 1, $false, "21-12=", 2) | Out-Null
$d.Content.Find.Execute("60+4=", $true, $false, $false, $false, $false, $true, 1, $false, "76-68=", 2) | Out-Null
$d.Content.Find.Execute("98-22=", $true, $false, $false, $false, $false, $true, 1, $false, "41+48=", 2) | Out-Null
$d.Content.Find.Execute("60-29=", $true, $false, $false, $false, $false, $true, 1, $false, "31+66=", 2) | Out-Null
$d.Content.Find.Execute("77+13=", $true, $false, $false, $false, $false, $true, 1, $false, "22+69=", 2) | Out-Null
$d.Content.Find.Execute("30-24=", $true, $false, $false, $false, $false, $true, 1, $false, "11+79=", 2) | Out-Null
$d.Content.Find.Execute("4+88=", $true, $false, $false, $false, $false, $true, 1, $false, "71-41=", 2) | Out-Null
$d.Content.Find.Execute("42+48=", $true, $false, $false, $false, $false, $true, 1, $false, "59-44=", 2) | Out-Null
$d.Content.Find.Execute("29+62=", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=", 2) | Out-Null
$d.Content.Find.Execute("44+22=", $true, $false, $false, $false, $false, $true, 1, $false, "20+48=", 2) | Out-Null
$d.Content.Find.Execute("52+3=", $true, $false, $false, $false, $false, $true, 1, $false, "24+18=", 2) | Out-Null
$d.Content.Find.Execute("23+67=", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=", 2) | Out-Null
$d.Content.Find.Execute("25+5=", $true, $false, $false, $false, $false, $true, 1, $false, "78-64=", 2) | Out-Null
$d.Content.Find.Execute("31+0=", $true, $false, $false, $false, $false, $true, 1, $false, "60+37=", 2) | Out-Null
$d.Content.Find.Execute("23+48=", $true, $false, $false, $false, $false, $true, 1, $false, "49+20=", 2) | Out-Null
$d.Content.Find.Execute("55+35=", $true, $false, $false, $false, $false, $true, 1, $false, "43+55=", 2) | Out-Null
$d.Content.Find.Execute("50+32=", $true, $false, $false, $false, $false, $true, 1, $false, "12+74=", 2) | Out-Null
$d.Content.Find.Execute("90-59=", $true, $false, $false, $false, $false, $true, 1, $false, "3+42=", 2) | Out-Null
$d.Content.Find.Execute("16+30=", $true, $false, $false, $false, $false, $true, 1, $false, "23-7=", 2) | Out-Null
$d.Content.Find.Execute("81-60=", $true, $false, $false, $false, $false, $true, 1, $false, "98-83=", 2) | Out-Null
